# This workbook uses one sheet to hold several different "record types"
# (feature / homolog / URL), each with its own mini-header row, stacked one
# after another. Only the "feature" block (rows 1-5) gains a new
# "featureType" column, inserted between "featureName" and "start". The
# other blocks are not column-shifted - we therefore move values cell by
# cell instead of doing a sheet-wide column insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Shift columns I..F <- H..E for rows 1-5 (the "feature" block), right
#     to left so we never clobber a value before it has been copied. ---
for ($r = 1; $r -le 5; $r++) {
    $ws.Cells.Item($r, 9).Value2 = $ws.Cells.Item($r, 8).Value2   # I <- H
    $ws.Cells.Item($r, 8).Value2 = $ws.Cells.Item($r, 7).Value2   # H <- G
    $ws.Cells.Item($r, 7).Value2 = $ws.Cells.Item($r, 6).Value2   # G <- F
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 5).Value2   # F <- E
}

# --- Fill the new "featureType" column (E). ---
$ws.Range("E1").Value2 = "featureType"
$ws.Range("E2").Value2 = "gene"
$ws.Range("E3").Value2 = "gene"
$ws.Range("E4").Value2 = "SNP"
$ws.Range("E5").Value2 = "SNP"

# Row 1 is a bold, centred header row (style carries through the whole
# row); make sure the new E1 cell picks up the same formatting as its
# neighbours.
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108   # xlCenter

# --- Row 9's trailing comment cell moves from H9 to I9 (this row's block
#     only has 3 data columns, so this is a simple relocation, not part of
#     a column insert). ---
$ws.Range("I9").Value2 = $ws.Range("H9").Value2
$ws.Range("I9").Font.Bold = $true
$ws.Range("H9").Clear()

# --- Column widths: shift the old E/F/G widths right onto F/G/H (mirroring
#     the column that conceptually got inserted), then size the new E and I
#     columns themselves. Read-then-write keeps every value consistent with
#     what this engine's own ColumnWidth getter reports. ---
$oldE = $ws.Range("E1").ColumnWidth
$oldF = $ws.Range("F1").ColumnWidth
$oldG = $ws.Range("G1").ColumnWidth
$ws.Range("H1").ColumnWidth = $oldG
$ws.Range("G1").ColumnWidth = $oldF
$ws.Range("F1").ColumnWidth = $oldE
$ws.Range("E1").ColumnWidth = 23.5
$ws.Range("I1").ColumnWidth = 52.16666666666667

# --- Update the saved cursor position recorded in the sheet view. ---
$ws.Range("F19").Select()
